$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Email" column (G) with header + mailto hyperlink value
$ws.Range("G1").Value = "Email"
$ws.Range("G2").Value = "adrianrentea01@gmail.com"

# Column width to match the other bestFit-style columns (~24.5 chars)
$ws.Columns("G").ColumnWidth = 23.667

# Add the hyperlink for the email address
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:adrianrentea01@gmail.com")

# The hyperlink write auto-applies Excel's built-in "Hyperlink" cell style;
# restore the cell to the plain/default style (no other cell carries an
# explicit style override here) and drop the now-unused named style so the
# workbook's style table stays as close as possible to the original.
$ws.Range("G2").Style = "Normal"
$wb.Styles("Hyperlink").Delete()

# Move the active selection to G11, matching the saved cursor position
$ws.Range("G11").Select() | Out-Null
